$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 55
$ws.Range("C5").Value = 54
$ws.Range("D5").Value = 48
